# Slide 1 title: "Automating " + "Geoprocessing" + " of Forestry Field Data"
#            ->  "Automating Forestry Field Data " + "Geoprocesses"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

$oldRun1 = "Automating "
$oldRun2 = "Geoprocessing"
$oldRun3 = " of Forestry Field Data"

$len1 = $oldRun1.Length
$len2 = $oldRun2.Length
$len3 = $oldRun3.Length

# Edit right-to-left so earlier character offsets stay valid while editing.

# Third run (" of Forestry Field Data") is removed entirely.
$run3 = $tr.Characters($len1 + $len2 + 1, $len3)
$run3.Text = ""

# Second run ("Geoprocessing" -> "Geoprocesses"), keeps its own run formatting (err="1").
$run2 = $tr.Characters($len1 + 1, $len2)
$run2.Text = "Geoprocesses"

# First run ("Automating " -> "Automating Forestry Field Data "), keeps its own run formatting.
$run1 = $tr.Characters(1, $len1)
$run1.Text = "Automating Forestry Field Data "
